# Contacts Template: add a "status" column and rename the first
# column header from "name" to "contact_name".
#
# Final header row (row 1):
#   A1 = contact_name
#   B1 = contact_number
#   C1 = image_name
#   D1 = status

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Activate()

$ws.Range("A1").Value = "contact_name"
$ws.Range("B1").Value = "contact_number"
$ws.Range("C1").Value = "image_name"
$ws.Range("D1").Value = "status"

$ws.Range("C5").Select()
